$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffix = "`n" + [char]8211 + ":" + [char]8211

$cells = @("H4", "H5", "H9", "H11")
foreach ($cellRef in $cells) {
    $cell = $ws.Range($cellRef)
    $current = $cell.Value()
    $cell.Value = $current + $suffix
}
